$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Huesca" (row 53) and "Huelva" (row 54) rows in full,
# since the shared-string table reordered the two labels while the
# worksheet cell data moved with them (net effect: rows 53 and 54 swap).
$row53 = $ws.Range("A53:E53").Value()
$row54 = $ws.Range("A54:E54").Value()

$ws.Range("A53:E53").Value = $row54
$ws.Range("A54:E54").Value = $row53

# Update the "last updated" timestamp string in A1.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 07:16"
